# Automatische test-sync: 2025-06-26 23:45:50
# Appends the newest test-mail log entry (row 44) to the "Logs" sheet
# and bumps the "Bestelling / Levering" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$row = 44
$logs.Cells.Item($row, 1).Value = "Wil je dit artikel voor me inkopen?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #12: Wil je dit artikel voor me inkopen?"
$logs.Cells.Item($row, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor je interesse in ons artikel. Helaas kan ik je op basis van dit bericht niet verder helpen. Kun je meer details geven over welk artikel je wilt inkopen en op welke manier? Zo kan ik je beter assisteren.`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-26 23:45:03"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 20

# The conditional-formatting rules were scoped to the old used range
# (…2:…43); extend each rule's AppliesTo range to cover the new row 44.
$cfColumns = @("D", "G", "H", "I")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "43")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "44")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
